# Fruta / hortaliza, semanal
# Insert a new weekly record at row 643 (pushing the existing rows 643:691
# down to 644:692) on the single data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 643:691 down to 644:692, leaving a blank row 643 to fill in.
$ws.Rows("643:643").Insert()

$ws.Range("A643").Value2 = 5
$ws.Range("B643").Value2 = "Macroferia Regional de Talca"
$ws.Range("C643").Value2 = "Maule"
$ws.Range("D643").Value2 = 45265
$ws.Range("E643").Value2 = 7
$ws.Range("F643").Value2 = 100112032
$ws.Range("G643").Value2 = "Zapallo italiano"
$ws.Range("H643").Value2 = "Sin especificar"
$ws.Range("I643").Value2 = "Primera"
$ws.Range("J643").Value2 = 400
$ws.Range("K643").Value2 = 6000
$ws.Range("L643").Value2 = 6000
$ws.Range("M643").Value2 = 6000
$ws.Range("N643").Value2 = "`$/caja 50 unidades"
$ws.Range("O643").Value2 = "Región del Maule"
$ws.Range("P643").Value2 = 120
$ws.Range("Q643").Value2 = 50
$ws.Range("R643").Value2 = "Hortaliza"
